$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "51.253.56"
$cell.ClearFormats()

$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.58%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.918.02"
$cell.ClearFormats()

$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.06%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.11%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "368.97"
$cell.ClearFormats()

$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.24%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "103.64"
$cell.ClearFormats()

$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.42%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.29%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.83%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "36.88"
$cell.ClearFormats()

$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.71%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.44%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.64%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "18.48"
$cell.ClearFormats()

$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.46%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.373.66"
$cell.ClearFormats()

$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.26%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.41"
$cell.ClearFormats()

$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.94%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.914.38"
$cell.ClearFormats()

$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.12%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.945"
$cell.ClearFormats()

$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.62%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "51.182.32"
$cell.ClearFormats()

$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.70%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.27"
$cell.ClearFormats()

$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -3.89%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.24"
$cell.ClearFormats()

$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.62%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.85"
$cell.ClearFormats()

$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -3.73%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0947"
$cell.ClearFormats()

$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.99%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "68.49"
$cell.ClearFormats()

$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.50%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "260.83"
$cell.ClearFormats()

$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.20%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.172"
$cell.ClearFormats()

$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.68%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "25.80"
$cell.ClearFormats()

$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.38%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.04"
$cell.ClearFormats()

$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -4.82%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.ClearFormats()

$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.56%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.09"
$cell.ClearFormats()

$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.18%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.44%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "34.81"
$cell.ClearFormats()

$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.17%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "50.87"
$cell.ClearFormats()

$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.33%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.24%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0424"
$cell.ClearFormats()

$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.05"
$cell.ClearFormats()

$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.39%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.68"
$cell.ClearFormats()

$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.92%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.15"
$cell.ClearFormats()

$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.65%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.86"
$cell.ClearFormats()

$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -5.04%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.17%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "22.26"
$cell.ClearFormats()

$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.77%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "119.87"
$cell.ClearFormats()

$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.51%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.08"
$cell.ClearFormats()

$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -3.69%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.026.23"
$cell.ClearFormats()

$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -3.16%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -5.87%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.18"
$cell.ClearFormats()

$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -3.85%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.194.63"
$cell.ClearFormats()

$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.55%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.239"
$cell.ClearFormats()

$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.20%  "
$cell.ClearFormats()

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0314"
$cell.ClearFormats()

$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -6.82%  "
$cell.ClearFormats()

